# This script updates the "想去人数" (F column) values on the "展览" and
# "全部类型" worksheets to reflect refreshed counts from a later data pull.

$wb = $excel.ActiveWorkbook

# Changes for the "展览" sheet (row -> new F value)
$exhibitionChanges = @{
    2  = 13751
    5  = 544
    8  = 1027
    9  = 13888
    10 = 14732
    12 = 3
    14 = 176
    20 = 18
    21 = 60
    22 = 16
    23 = 1144
    26 = 5704
    28 = 1055
    29 = 5403
    31 = 48
    32 = 248
}

# Changes for the "全部类型" sheet (row -> new F value)
$allTypesChanges = @{
    2  = 13751
    6  = 544
    9  = 1027
    10 = 13888
    11 = 14732
    13 = 3
    15 = 176
    21 = 18
    22 = 60
    23 = 16
    24 = 1144
    27 = 5704
    29 = 1055
    30 = 5403
    32 = 48
    33 = 248
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionChanges.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionChanges[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesChanges.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesChanges[$row]
}
